$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tmp = $ws.Range("Z1")
$tmp.Formula = "=`"27.239.92`""
$tmp.Copy()
$ws.Range("D2").PasteSpecial(-4163)
$tmp.Formula = "=`"  -0.09%  `""
$tmp.Copy()
$ws.Range("E2").PasteSpecial(-4163)
$tmp.Formula = "=`"1.902.43`""
$tmp.Copy()
$ws.Range("D3").PasteSpecial(-4163)
$tmp.Formula = "=`"  +0.18%  `""
$tmp.Copy()
$ws.Range("E3").PasteSpecial(-4163)
$tmp.Formula = "=`"  -0.06%  `""
$tmp.Copy()
$ws.Range("E4").PasteSpecial(-4163)
$tmp.Formula = "=`"306.34`""
$tmp.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$tmp.Formula = "=`"  -0.51%  `""
$tmp.Copy()
$ws.Range("E5").PasteSpecial(-4163)
$tmp.Formula = "=`"  -0.05%  `""
$tmp.Copy()
$ws.Range("E6").PasteSpecial(-4163)
$tmp.Formula = "=`"0.5350`""
$tmp.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$tmp.Formula = "=`"  +2.63%  `""
$tmp.Copy()
$ws.Range("E7").PasteSpecial(-4163)
$tmp.Formula = "=`"  +0.89%  `""
$tmp.Copy()
$ws.Range("E8").PasteSpecial(-4163)
$tmp.Formula = "=`"0.07285`""
$tmp.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$tmp.Formula = "=`"  -0.01%  `""
$tmp.Copy()
$ws.Range("E9").PasteSpecial(-4163)
$tmp.Formula = "=`"22.21`""
$tmp.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$tmp.Formula = "=`"  +4.67%  `""
$tmp.Copy()
$ws.Range("E10").PasteSpecial(-4163)
$tmp.Formula = "=`"0.9022`""
$tmp.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$tmp.Formula = "=`"  +0.04%  `""
$tmp.Copy()
$ws.Range("E11").PasteSpecial(-4163)
$tmp.Formula = "=`"0.08217`""
$tmp.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$tmp.Formula = "=`"  +0.20%  `""
$tmp.Copy()
$ws.Range("E12").PasteSpecial(-4163)
$tmp.Formula = "=`"96.09`""
$tmp.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$tmp.Formula = "=`"  -0.77%  `""
$tmp.Copy()
$ws.Range("E13").PasteSpecial(-4163)
$tmp.Formula = "=`"5.333`""
$tmp.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$tmp.Formula = "=`"  +0.67%  `""
$tmp.Copy()
$ws.Range("E14").PasteSpecial(-4163)
$tmp.Formula = "=`"  -0.11%  `""
$tmp.Copy()
$ws.Range("E15").PasteSpecial(-4163)
$tmp.Formula = "=`"14.84`""
$tmp.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$tmp.Formula = "=`"  +1.84%  `""
$tmp.Copy()
$ws.Range("E16").PasteSpecial(-4163)
$tmp.Formula = "=`"0.000008649`""
$tmp.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$tmp.Formula = "=`"  +0.50%  `""
$tmp.Copy()
$ws.Range("E17").PasteSpecial(-4163)
$tmp.Formula = "=`"  -0.05%  `""
$tmp.Copy()
$ws.Range("E18").PasteSpecial(-4163)
$tmp.Formula = "=`"27.261.54`""
$tmp.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$tmp.Formula = "=`"  -0.13%  `""
$tmp.Copy()
$ws.Range("E19").PasteSpecial(-4163)
$tmp.Formula = "=`"5.035`""
$tmp.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$tmp.Formula = "=`"  -1.19%  `""
$tmp.Copy()
$ws.Range("E20").PasteSpecial(-4163)
$tmp.Formula = "=`"1.086.24`""
$tmp.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$tmp.Formula = "=`"  -42.99%  `""
$tmp.Copy()
$ws.Range("E21").PasteSpecial(-4163)
$tmp.Formula = "=`"  +0.48%  `""
$tmp.Copy()
$ws.Range("E22").PasteSpecial(-4163)
$tmp.Formula = "=`"6.494`""
$tmp.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$tmp.Formula = "=`"  +1.23%  `""
$tmp.Copy()
$ws.Range("E23").PasteSpecial(-4163)
$tmp.Formula = "=`"149.60`""
$tmp.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$tmp.Formula = "=`"  +1.53%  `""
$tmp.Copy()
$ws.Range("E24").PasteSpecial(-4163)
$tmp.Formula = "=`"2.294`""
$tmp.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$tmp.Formula = "=`"  -0.45%  `""
$tmp.Copy()
$ws.Range("E25").PasteSpecial(-4163)
$tmp.Formula = "=`"18.35`""
$tmp.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$tmp.Formula = "=`"  +0.59%  `""
$tmp.Copy()
$ws.Range("E26").PasteSpecial(-4163)
$tmp.Formula = "=`"1.747`""
$tmp.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$tmp.Formula = "=`"  -0.01%  `""
$tmp.Copy()
$ws.Range("E27").PasteSpecial(-4163)
$tmp.Formula = "=`"116.71`""
$tmp.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$tmp.Formula = "=`"  +0.98%  `""
$tmp.Copy()
$ws.Range("E28").PasteSpecial(-4163)
$tmp.Formula = "=`"4.808`""
$tmp.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$tmp.Formula = "=`"  -0.50%  `""
$tmp.Copy()
$ws.Range("E29").PasteSpecial(-4163)
$tmp.Formula = "=`"4.769`""
$tmp.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$tmp.Formula = "=`"  -3.04%  `""
$tmp.Copy()
$ws.Range("E30").PasteSpecial(-4163)
$tmp.Formula = "=`"0.09219`""
$tmp.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$tmp.Formula = "=`"  -0.28%  `""
$tmp.Copy()
$ws.Range("E31").PasteSpecial(-4163)
$tmp.Formula = "=`"0.8266`""
$tmp.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$tmp.Formula = "=`"  +3.69%  `""
$tmp.Copy()
$ws.Range("E32").PasteSpecial(-4163)
$tmp.Formula = "=`"0.05062`""
$tmp.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$tmp.Formula = "=`"  +0.01%  `""
$tmp.Copy()
$ws.Range("E33").PasteSpecial(-4163)
$tmp.Formula = "=`"1.216`""
$tmp.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$tmp.Formula = "=`"  -1.31%  `""
$tmp.Copy()
$ws.Range("E34").PasteSpecial(-4163)
$tmp.Formula = "=`"2.992`""
$tmp.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$tmp.Formula = "=`"  +0.84%  `""
$tmp.Copy()
$ws.Range("E35").PasteSpecial(-4163)
$tmp.Formula = "=`"3.344`""
$tmp.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$tmp.Formula = "=`"  -2.67%  `""
$tmp.Copy()
$ws.Range("E36").PasteSpecial(-4163)
$tmp.Formula = "=`"2.674`""
$tmp.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$tmp.Formula = "=`"  +3.10%  `""
$tmp.Copy()
$ws.Range("E37").PasteSpecial(-4163)
$tmp.Formula = "=`"0.5742`""
$tmp.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$tmp.Formula = "=`"  +0.96%  `""
$tmp.Copy()
$ws.Range("E38").PasteSpecial(-4163)
$tmp.Formula = "=`"  +0.42%  `""
$tmp.Copy()
$ws.Range("E39").PasteSpecial(-4163)
$tmp.Formula = "=`"  -0.10%  `""
$tmp.Copy()
$ws.Range("E40").PasteSpecial(-4163)
$tmp.Formula = "=`"9.370`""
$tmp.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$tmp.Formula = "=`"  +4.31%  `""
$tmp.Copy()
$ws.Range("E41").PasteSpecial(-4163)
$tmp.Formula = "=`"6.589`""
$tmp.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$tmp.Formula = "=`"  +0.26%  `""
$tmp.Copy()
$ws.Range("E42").PasteSpecial(-4163)
$tmp.Formula = "=`"116.95`""
$tmp.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$tmp.Formula = "=`"  +1.40%  `""
$tmp.Copy()
$ws.Range("E43").PasteSpecial(-4163)
$tmp.Formula = "=`"0.1524`""
$tmp.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$tmp.Formula = "=`"  +0.41%  `""
$tmp.Copy()
$ws.Range("E44").PasteSpecial(-4163)
$tmp.Formula = "=`"0.4954`""
$tmp.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$tmp.Formula = "=`"  +1.36%  `""
$tmp.Copy()
$ws.Range("E45").PasteSpecial(-4163)
$tmp.Formula = "=`"  -0.09%  `""
$tmp.Copy()
$ws.Range("E46").PasteSpecial(-4163)
$tmp.Formula = "=`"  +0.39%  `""
$tmp.Copy()
$ws.Range("E47").PasteSpecial(-4163)
$tmp.Formula = "=`"1.638`""
$tmp.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$tmp.Formula = "=`"  +0.83%  `""
$tmp.Copy()
$ws.Range("E48").PasteSpecial(-4163)
$tmp.Formula = "=`"38.30`""
$tmp.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$tmp.Formula = "=`"  +0.36%  `""
$tmp.Copy()
$ws.Range("E49").PasteSpecial(-4163)
$tmp.Formula = "=`"  +3.88%  `""
$tmp.Copy()
$ws.Range("E50").PasteSpecial(-4163)
$tmp.Formula = "=`"63.24`""
$tmp.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$tmp.Formula = "=`"  -0.69%  `""
$tmp.Copy()
$ws.Range("E51").PasteSpecial(-4163)
$tmp.ClearContents()
$excel.CutCopyMode = 0
